# Applies the "resolving issue #1" edit to Metrics.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the "Enriched Entropy Enriched Sequence entropy" label to
#    "Enriched  Sequence entropy"
$ws.Range("D19").Value = "Enriched  Sequence entropy"

# 2. Drop the resolved review notes in the "Notes" column (G11, G12, G30)
#    together with the red "needs attention" highlight that was on the
#    corresponding "Metrics implemented in tool" cells (E11, E12, E30).
$ws.Range("G11").Clear()
$ws.Range("G12").Clear()
$ws.Range("G30").Clear()

$ws.Range("E9").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Those rows no longer need the extra height that was reserved for the
#    wrapped note text, so let them size back to the sheet's default.
$ws.Rows.Item(11).EntireRow.AutoFit()
$ws.Rows.Item(12).EntireRow.AutoFit()
$ws.Rows.Item(30).EntireRow.AutoFit()

# 4. Update the saved view position/selection.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C22").Select() | Out-Null
